$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1265
$ws.Range("F5").Value = 2048
$ws.Range("F6").Value = 104
$ws.Range("F7").Value = 790
$ws.Range("F9").Value = 650
$ws.Range("F14").Value = 31
$ws.Range("F17").Value = 297
$ws.Range("F18").Value = 6
$ws.Range("F24").Value = 623
$ws.Range("F28").Value = 4915
$ws.Range("F31").Value = 2405
$ws.Range("F32").Value = 5756
$ws.Range("F40").Value = 642
$ws.Range("F47").Value = 300

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 116
$ws.Range("F41").Value = 887

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 735
$ws.Range("F7").Value = 347

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1265
$ws.Range("F7").Value = 735
$ws.Range("F9").Value = 347
$ws.Range("F13").Value = 2048
$ws.Range("F15").Value = 790
$ws.Range("F22").Value = 31
$ws.Range("F25").Value = 300
$ws.Range("F27").Value = 116
$ws.Range("F30").Value = 623
$ws.Range("F36").Value = 4915
$ws.Range("F37").Value = 2405
$ws.Range("F38").Value = 5756
$ws.Range("F44").Value = 642
$ws.Range("F48").Value = 887
